# Update loading_percent results for Case_2_10 (380 kV case) - rows 2:25, cols B:O.
# Columns D, K, M, O remain 0 (unused line slots); all other columns get
# refreshed load-flow results from the re-run power-flow case.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$values = New-Object 'object[,]' 24,14
$values[0,0] = 18.34672923810091
$values[0,1] = 9.184794409951568
$values[0,2] = 0
$values[0,3] = 14.42623828828769
$values[0,4] = 38.76547630129511
$values[0,5] = 38.01008693309698
$values[0,6] = 16.37295692481242
$values[0,7] = 25.16386250524436
$values[0,8] = 7.900487943053969
$values[0,9] = 0
$values[0,10] = 12.82890023916018
$values[0,11] = 0
$values[0,12] = 18.03691089615791
$values[0,13] = 0
$values[1,0] = 17.88285286626935
$values[1,1] = 8.820871057815937
$values[1,2] = 0
$values[1,3] = 14.44510135219539
$values[1,4] = 38.73927107487006
$values[1,5] = 37.88981539016788
$values[1,6] = 16.41104399871963
$values[1,7] = 25.25057699658393
$values[1,8] = 7.914883362388222
$values[1,9] = 0
$values[1,10] = 12.80832424424728
$values[1,11] = 0
$values[1,12] = 18.10086277752859
$values[1,13] = 0
$values[2,0] = 17.59537166305229
$values[2,1] = 8.587701125750435
$values[2,2] = 0
$values[2,3] = 14.45793278204563
$values[2,4] = 38.73332471907418
$values[2,5] = 37.83003543434909
$values[2,6] = 16.43805492014402
$values[2,7] = 25.30978788258752
$values[2,8] = 7.924150144783411
$values[2,9] = 0
$values[2,10] = 12.79749856521316
$values[2,11] = 0
$values[2,12] = 18.14202311622726
$values[2,13] = 0
$values[3,0] = 17.47773128296563
$values[3,1] = 8.490314984490217
$values[3,2] = 0
$values[3,3] = 14.46347649867637
$values[3,4] = 38.7334535712024
$values[3,5] = 37.80922410993161
$values[3,6] = 16.4499711931204
$values[3,7] = 25.33541227079636
$values[3,8] = 7.928034422362884
$values[3,9] = 0
$values[3,10] = 12.79354487962372
$values[3,11] = 0
$values[3,12] = 18.15927388707626
$values[3,13] = 0
$values[4,0] = 17.45817301962876
$values[4,1] = 8.474003541752536
$values[4,2] = 0
$values[4,3] = 14.46441606019939
$values[4,4] = 38.73362910403995
$values[4,5] = 37.80598303580231
$values[4,6] = 16.45200471220204
$values[4,7] = 25.33975734962036
$values[4,8] = 7.928685937104138
$values[4,9] = 0
$values[4,10] = 12.79291611567097
$values[4,11] = 0
$values[4,12] = 18.16216725245175
$values[4,13] = 0
$values[5,0] = 17.5937868628477
$values[5,1] = 8.58639721926459
$values[5,2] = 0
$values[5,3] = 14.45800627107965
$values[5,4] = 38.7333161234881
$values[5,5] = 37.82974038182216
$values[5,6] = 16.4382119497246
$values[5,7] = 25.31012741394937
$values[5,8] = 7.924202091771967
$values[5,9] = 0
$values[5,10] = 12.79744338649918
$values[5,11] = 0
$values[5,12] = 18.14225383049313
$values[5,13] = 0
$values[6,0] = 18.1874415499258
$values[6,1] = 9.061373681501353
$values[6,2] = 0
$values[6,3] = 14.43248333002601
$values[6,4] = 38.75433675924572
$values[6,5] = 37.96570657635471
$values[6,6] = 16.38533545445738
$values[6,7] = 25.19251995307339
$values[6,8] = 7.905362895135101
$values[6,9] = 0
$values[6,10] = 12.82143237827346
$values[6,11] = 0
$values[6,12] = 18.05856934564001
$values[6,13] = 0
$values[7,0] = 19.3230638044944
$values[7,1] = 9.912823498321565
$values[7,2] = 0
$values[7,3] = 14.39231914356174
$values[7,4] = 38.8759045701733
$values[7,5] = 38.34316495663327
$values[7,6] = 16.3105161892921
$values[7,7] = 25.00947168366596
$values[7,8] = 7.871797275755632
$values[7,9] = 0
$values[7,10] = 12.88266953927435
$values[7,11] = 0
$values[7,12] = 17.90942362338756
$values[7,13] = 0
$values[8,0] = 20.13074574548422
$values[8,1] = 10.48643288421343
$values[8,2] = 0
$values[8,3] = 14.36880047586209
$values[8,4] = 39.01389354994728
$values[8,5] = 38.68663156718811
$values[8,6] = 16.27328345111251
$values[8,7] = 24.90429606007138
$values[8,8] = 7.849171228402042
$values[8,9] = 0
$values[8,10] = 12.93609646106281
$values[8,11] = 0
$values[8,12] = 17.80887314947211
$values[8,13] = 0
$values[9,0] = 20.4906432456987
$values[9,1] = 10.73558376267644
$values[9,2] = 0
$values[9,3] = 14.35939389311689
$values[9,4] = 39.08712947929475
$values[9,5] = 38.85683666652152
$values[9,6] = 16.26022208459126
$values[9,7] = 24.86287552854487
$values[9,8] = 7.839314567303017
$values[9,9] = 0
$values[9,10] = 12.96218171446583
$values[9,11] = 0
$values[9,12] = 17.76507063334423
$values[9,13] = 0
$values[10,0] = 20.62571465974166
$values[10,1] = 10.82820104369992
$values[10,2] = 0
$values[10,3] = 14.35601700187968
$values[10,4] = 39.11635421858168
$values[10,5] = 38.92325439319618
$values[10,6] = 16.25583516578131
$values[10,7] = 24.84811932780274
$values[10,8] = 7.835644408457503
$values[10,9] = 0
$values[10,10] = 12.97231078683247
$values[10,11] = 0
$values[10,12] = 17.74876105753753
$values[10,13] = 0
$values[11,0] = 20.59668059659154
$values[11,1] = 10.80833174034816
$values[11,2] = 0
$values[11,3] = 14.35673605025593
$values[11,4] = 39.10999402208147
$values[11,5] = 38.90886349329988
$values[11,6] = 16.25675507357304
$values[11,7] = 24.85125595039596
$values[11,8] = 7.836432075700718
$values[11,9] = 0
$values[11,10] = 12.97011822150435
$values[11,11] = 0
$values[11,12] = 17.75226129495696
$values[11,13] = 0
$values[12,0] = 20.50178070112714
$values[12,1] = 10.74323831373444
$values[12,2] = 0
$values[12,3] = 14.3591123666899
$values[12,4] = 39.08950399735516
$values[12,5] = 38.86226177730907
$values[12,6] = 16.25984995482953
$values[12,7] = 24.86164288321837
$values[12,8] = 7.839011373684655
$values[12,9] = 0
$values[12,10] = 12.96301003519895
$values[12,11] = 0
$values[12,12] = 17.76372328059404
$values[12,13] = 0
$values[13,0] = 20.44348993963275
$values[13,1] = 10.70314035091275
$values[13,2] = 0
$values[13,3] = 14.36059202563849
$values[13,4] = 39.07714713717675
$values[13,5] = 38.83397140285779
$values[13,6] = 16.26181852398735
$values[13,7] = 24.86812629009869
$values[13,8] = 7.840599376730157
$values[13,9] = 0
$values[13,10] = 12.95868862171899
$values[13,11] = 0
$values[13,12] = 17.77078017782339
$values[13,13] = 0
$values[14,0] = 20.10706268065329
$values[14,1] = 10.46990991194482
$values[14,2] = 0
$values[14,3] = 14.36944116595553
$values[14,4] = 39.0093169223475
$values[14,5] = 38.67578563422786
$values[14,6] = 16.27421519372913
$values[14,7] = 24.90713269919253
$values[14,8] = 7.849824128056643
$values[14,9] = 0
$values[14,10] = 12.93442714977913
$values[14,11] = 0
$values[14,12] = 17.81177464489421
$values[14,13] = 0
$values[15,0] = 19.89865138373856
$values[15,1] = 10.32378433470521
$values[15,2] = 0
$values[15,3] = 14.37520033536476
$values[15,4] = 38.97037689099187
$values[15,5] = 38.58229044688179
$values[15,6] = 16.28281421621353
$values[15,7] = 24.93271095100502
$values[15,8] = 7.855594639973345
$values[15,9] = 0
$values[15,10] = 12.91999624940104
$values[15,11] = 0
$values[15,12] = 17.83741899532503
$values[15,13] = 0
$values[16,0] = 19.77807801175672
$values[16,1] = 10.23862907261272
$values[16,2] = 0
$values[16,3] = 14.3786345218622
$values[16,4] = 38.94896543904503
$values[16,5] = 38.52983136908085
$values[16,6] = 16.28812488084781
$values[16,7] = 24.9480272889157
$values[16,8] = 7.858954749295574
$values[16,9] = 0
$values[16,10] = 12.91186393433484
$values[16,11] = 0
$values[16,12] = 17.85235147226679
$values[16,13] = 0
$values[17,0] = 19.73713780953115
$values[17,1] = 10.2096079032654
$values[17,2] = 0
$values[17,3] = 14.37981819360222
$values[17,4] = 38.94188556340865
$values[17,5] = 38.51229702907379
$values[17,6] = 16.28998556649112
$values[17,7] = 24.9533167882335
$values[17,8] = 7.860099488202554
$values[17,9] = 0
$values[17,10] = 12.90913946602977
$values[17,11] = 0
$values[17,12] = 17.85743873973174
$values[17,13] = 0
$values[18,0] = 19.92091063695239
$values[18,1] = 10.33945457185841
$values[18,2] = 0
$values[18,3] = 14.37457467470682
$values[18,4] = 38.97442018117984
$values[18,5] = 38.59210717189964
$values[18,6] = 16.28186107335978
$values[18,7] = 24.92992551174522
$values[18,8] = 7.854976111797218
$values[18,9] = 0
$values[18,10] = 12.92151509559679
$values[18,11] = 0
$values[18,12] = 17.83467022892114
$values[18,13] = 0
$values[19,0] = 20.52968904806928
$values[19,1] = 10.76240506012894
$values[19,2] = 0
$values[19,3] = 14.35840936440977
$values[19,4] = 39.09548203736428
$values[19,5] = 38.87589687131374
$values[19,6] = 16.25892572532261
$values[19,7] = 24.85856674120441
$values[19,8] = 7.838252082163004
$values[19,9] = 0
$values[19,10] = 12.965091106419
$values[19,11] = 0
$values[19,12] = 17.76034909655499
$values[19,13] = 0
$values[20,0] = 20.92043136185727
$values[20,1] = 11.02872680139616
$values[20,2] = 0
$values[20,3] = 14.34892345301533
$values[20,4] = 39.18329141794293
$values[20,5] = 39.07279895099852
$values[20,6] = 16.24719589353442
$values[20,7] = 24.81734593902067
$values[20,8] = 7.827685190615963
$values[20,9] = 0
$values[20,10] = 12.9950319811379
$values[20,11] = 0
$values[20,12] = 17.71339277592345
$values[20,13] = 0
$values[21,0] = 20.71257810353256
$values[21,1] = 10.88752048926965
$values[21,2] = 0
$values[21,3] = 14.35388774191898
$values[21,4] = 39.13563566536247
$values[21,5] = 38.96667811820323
$values[21,6] = 16.25315753051975
$values[21,7] = 24.8388490534321
$values[21,8] = 7.833291821698944
$values[21,9] = 0
$values[21,10] = 12.97891995658083
$values[21,11] = 0
$values[21,12] = 17.73830671629548
$values[21,13] = 0
$values[22,0] = 19.9108495698174
$values[22,1] = 10.33237362366282
$values[22,2] = 0
$values[22,3] = 14.37485715222175
$values[22,4] = 38.97258916944946
$values[22,5] = 38.58766500101427
$values[22,6] = 16.28229084621921
$values[22,7] = 24.93118290571291
$values[22,8] = 7.855255615754304
$values[22,9] = 0
$values[22,10] = 12.92082791339311
$values[22,11] = 0
$values[22,12] = 17.83591235673286
$values[22,13] = 0
$values[23,0] = 19.01989955967448
$values[23,1] = 9.691407384150132
$values[23,2] = 0
$values[23,3] = 14.40213001167409
$values[23,4] = 38.83443934552366
$values[23,5] = 38.229312910819
$values[23,6] = 16.32765108481728
$values[23,7] = 25.05386537291409
$values[23,8] = 7.880518588303784
$values[23,9] = 0
$values[23,10] = 12.86460497485872
$values[23,11] = 0
$values[23,12] = 17.94817965195017
$values[23,13] = 0
$ws.Range("B2:O25").Value = $values
Write-Host "Done"